# Natmi following Dr Hou advice
# Update the Bmp6-Acvr2a LR-pair table: recompute the cell-cell communication
# scores now that the "ECs" sending cluster is included (3 sending clusters x
# 4 target clusters = 12 data rows).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
  @("ECs", "Bmp6", "Acvr2a", "ECs", 2, 0.6666666666666666, 19.72450833333333, 59.173525, 0.5834853563809828, 0.5834853563809829, 3, 1, 10.793878, 32.381634, 0.166691744666625, 0.166691744666625, 212.9039365599833, 1916.13542903985, 0.09726219204257344, 0.09726219204257346),
  @("ECs", "Bmp6", "Acvr2a", "FAPs", 2, 0.6666666666666666, 19.72450833333333, 59.173525, 0.5834853563809828, 0.5834853563809829, 3, 1, 29.04277733333333, 87.128332, 0.4485126868821051, 0.4485126868821051, 572.8545035344778, 5155.6905318103, 0.2617005849467972, 0.2617005849467973),
  @("ECs", "Bmp6", "Acvr2a", "M2", 2, 0.6666666666666666, 19.72450833333333, 59.173525, 0.5834853563809828, 0.5834853563809829, 3, 1, 5.345855666666666, 16.037567, 0.08255698348755007, 0.08255698348755008, 105.4443746459639, 948.9993718136749, 0.04817079093197206, 0.04817079093197208),
  @("ECs", "Bmp6", "Acvr2a", "sCs", 2, 0.6666666666666666, 19.72450833333333, 59.173525, 0.5834853563809828, 0.5834853563809829, 3, 1, 19.57101366666667, 58.713041, 0.3022385849637199, 0.3022385849637199, 386.0286221599472, 3474.257599439525, 0.1763517884596401, 0.1763517884596401),
  @("FAPs", "Bmp6", "Acvr2a", "ECs", 3, 1, 0.7684289999999999, 2.305287, 0.02273146997336134, 0.02273146997336134, 3, 1, 10.793878, 32.381634, 0.166691744666625, 0.166691744666625, 8.294328877661998, 74.648959898958, 0.003789148388696601, 0.003789148388696601),
  @("FAPs", "Bmp6", "Acvr2a", "FAPs", 3, 1, 0.7684289999999999, 2.305287, 0.02273146997336134, 0.02273146997336134, 3, 1, 29.04277733333333, 87.128332, 0.4485126868821051, 0.4485126868821051, 22.317312343476, 200.855811091284, 0.01019535267453219, 0.01019535267453219),
  @("FAPs", "Bmp6", "Acvr2a", "M2", 3, 1, 0.7684289999999999, 2.305287, 0.02273146997336134, 0.02273146997336134, 3, 1, 5.345855666666666, 16.037567, 0.08255698348755007, 0.08255698348755008, 4.107910524080999, 36.97119471672899, 0.001876641591238532, 0.001876641591238533),
  @("FAPs", "Bmp6", "Acvr2a", "sCs", 3, 1, 0.7684289999999999, 2.305287, 0.02273146997336134, 0.02273146997336134, 3, 1, 19.57101366666667, 58.713041, 0.3022385849637199, 0.3022385849637199, 15.038934460863, 135.350410147767, 0.00687032731889402, 0.00687032731889402),
  @("sCs", "Bmp6", "Acvr2a", "ECs", 3, 1, 13.31169566666667, 39.935087, 0.3937831736456558, 0.3937831736456558, 3, 1, 10.793878, 32.381634, 0.166691744666625, 0.166691744666625, 143.6848189991287, 1293.163370992158, 0.06564040423535489, 0.0656404042353549),
  @("sCs", "Bmp6", "Acvr2a", "FAPs", 3, 1, 13.31169566666667, 39.935087, 0.3937831736456558, 0.3937831736456558, 3, 1, 29.04277733333333, 87.128332, 0.4485126868821051, 0.4485126868821051, 386.6086131760982, 3479.477518584884, 0.1766167492607756, 0.1766167492607757),
  @("sCs", "Bmp6", "Acvr2a", "M2", 3, 1, 13.31169566666667, 39.935087, 0.3937831736456558, 0.3937831736456558, 3, 1, 5.345855666666666, 16.037567, 0.08255698348755007, 0.08255698348755008, 71.16240371259211, 640.461633413329, 0.03250955096433947, 0.03250955096433947),
  @("sCs", "Bmp6", "Acvr2a", "sCs", 3, 1, 13.31169566666667, 39.935087, 0.3937831736456558, 0.3937831736456558, 3, 1, 19.57101366666667, 58.713041, 0.3022385849637199, 0.3022385849637199, 260.5233778188408, 2344.710400369567, 0.1190164691851858, 0.1190164691851858)
)

$startRow = 2
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $rowValues = $rows[$i]
    for ($c = 1; $c -le $rowValues.Count; $c++) {
        $ws.Cells.Item($r, $c).Value = $rowValues[$c - 1]
    }
}
